$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.013.75"
$ws.Range("E2").Value = "  +2.57%  "
$ws.Range("D3").Value = "2.231.85"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'293.20"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'86.59"
$ws.Range("E6").Value = "  +5.72%  "
$ws.Range("D7").Value = "'0.514"
$ws.Range("E7").Value = "  +1.18%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "'30.67"
$ws.Range("E10").Value = "  +6.61%  "
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").Value = "'47.00"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").Value = "'6.39"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "2.575.04"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "'14.06"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "2.216.29"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "'0.727"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").Value = "39.927.83"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("E21").Value = "  +9.19%  "
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").Value = "'65.18"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").Value = "'234.88"
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'2.46"
$ws.Range("E26").Value = "  +3.49%  "
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  +5.29%  "
$ws.Range("D28").Value = "'22.76"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("D31").Value = "'33.22"
$ws.Range("E31").Value = "  +5.32%  "
$ws.Range("D32").Value = "'152.39"
$ws.Range("E32").Value = "  +3.43%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").Value = "'0.0720"
$ws.Range("E36").Value = "  +3.25%  "
$ws.Range("D37").Value = "'16.16"
$ws.Range("E37").Value = "  +10.06%  "
$ws.Range("E38").Value = "  +5.19%  "
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").Value = "'0.0999"
$ws.Range("E40").Value = "  +4.98%  "
$ws.Range("E41").Value = "  +5.90%  "
$ws.Range("E42").Value = "  +4.34%  "
$ws.Range("D43").Value = "2.039.60"
$ws.Range("E43").Value = "  +7.00%  "
$ws.Range("E44").Value = "  +7.77%  "
$ws.Range("D45").Value = "'0.0270"
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("D46").Value = "'10.04"
$ws.Range("E46").Value = "  +11.65%  "
$ws.Range("D47").Value = "'16.87"
$ws.Range("E47").Value = "  +6.83%  "
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").Value = "2.463.15"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "'71.12"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").Value = "'89.22"
$ws.Range("E51").Value = "  +3.20%  "
